# feat: add works and third course test scripts.
# Append two new benchmark rows (14-15) below the first results table on
# Sheet1, and leave the selection where the author ended up (H20).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 14: label only ("base") - new work entry, no timing data yet.
$ws.Range("B14").Value = "base"

# Row 15: "attention work3" entry with a controlnet timing value.
$ws.Range("B15").Value = "attention work3"
$ws.Range("C15").Value = 8.74451

# Match the author's final cell selection.
$ws.Range("H20").Select() | Out-Null
